$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1): relabel the report columns
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Product Name"
$ws.Range("B1").Value = "Category Name"
$ws.Range("C1").Value = "From"
$ws.Range("D1").Value = "To"
$ws.Range("E1").Value = "Total Quantity Sold"

# ---------------------------------------------------------------------------
# Data row (row 2): new sample data
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Ninja dildo"
$ws.Range("B2").Value = "dildos"
$ws.Range("C2").Value = "8/3/2014 12:00:00 AM"

# D2 must stay text ("666") rather than become a number, so force a text
# number format before assigning a purely-numeric-looking string.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "666"

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.5
$ws.Columns.Item(2).ColumnWidth = 14.333
$ws.Columns.Item(3).ColumnWidth = 20.667
$ws.Columns.Item(4).ColumnWidth = 3.833
$ws.Columns.Item(5).ColumnWidth = 17.833

# ---------------------------------------------------------------------------
# Header / footer: move title to the right header, swap footer pieces
# ---------------------------------------------------------------------------
$ws.PageSetup.CenterHeader = ""
$ws.PageSetup.RightHeader = "&24&U&""Arial,Regular Bold"" Sales Report"

$ws.PageSetup.LeftFooter = "Page &P of &N"
$ws.PageSetup.CenterFooter = "&A"
$ws.PageSetup.RightFooter = ""
